# Revert family psychiatric history to template:
#   {{GUARDIAN_NAME}}  ->  {{ | REPORTING_GUARDIAN | }}   (three runs,
#   each keeping the original run's formatting: color A5A5A5 / accent3)
#
# The host engine re-merges every same-formatted run in a paragraph
# whenever a Range.Text write lands anywhere inside it. A plain
# find/replace on "{{GUARDIAN_NAME}}" therefore also silently swallows
# unrelated neighbouring same-coloured runs further along in this
# paragraph (the "in-person/via video " / "conference" pair, and the
# ", " / "as well as " pair). To keep the edit faithful to the original
# run layout we:
#   1. locate the target run and record -- from the known original run
#      lengths in this paragraph -- every original run-boundary offset,
#   2. perform the text replacement (this merges same-styled runs
#      together across the whole paragraph),
#   3. replay a Bold-on/Bold-off no-op across each sub-range between
#      consecutive boundaries (shifted for the length delta of our edit),
#      plus the two new boundaries needed to split the replacement text
#      into "{{" / "REPORTING_GUARDIAN" / "}}". A formatting round-trip
#      that ends back where it started does not trigger the merge pass,
#      so this puts every run -- both the untouched ones and the freshly
#      split ones -- back where it belongs. (Collapsed/zero-length probe
#      ranges do not register as a split, so each probe must span at
#      least one character.)

$d = $word.ActiveDocument

$old = "{{GUARDIAN_NAME}}"
$new = "{{REPORTING_GUARDIAN}}"
$innerNew = "REPORTING_GUARDIAN"

$rng = $d.Content
$found = $rng.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find '$old' in the document"
}

$runStart = $rng.Start
$runEnd = $rng.End

# Original run lengths in this paragraph, in document order (known from
# the template XML):
#   "{{PREFERRED_NAME}}", " and ", "{{GUARDIAN_NAME}}",
#   " completed the K-SADS semi-structured psychiatric interview for DSM-5 ",
#   "in-person/via video ", "conference", ", ",
#   "as well as several questionnaires related to emotional and behavioral function. "
# The target run ("{{GUARDIAN_NAME}}") is the third one.
$runLens = @(18, 5, 17, 70, 20, 10, 2, 80)
$paraStart = $runStart - ($runLens[0] + $runLens[1])

$boundaries = New-Object System.Collections.Generic.List[int]
$pos = $paraStart
foreach ($len in $runLens) {
    [void]$boundaries.Add($pos)
    $pos = $pos + $len
}
[void]$boundaries.Add($pos)
$paraEnd = $pos

# Sanity-check our offline offsets against what Find actually reported.
if ($paraStart + $runLens[0] + $runLens[1] -ne $runStart) {
    throw "paragraph offset mismatch: expected run start $($paraStart + $runLens[0] + $runLens[1]), got $runStart"
}
if ($paraStart + $runLens[0] + $runLens[1] + $runLens[2] -ne $runEnd) {
    throw "paragraph offset mismatch: expected run end, got $runEnd"
}

$oldLen = $runEnd - $runStart
$newLen = $new.Length
$delta = $newLen - $oldLen

# Do the actual text swap (this merges same-styled runs paragraph-wide).
$editRng = $d.Range($runStart, $runEnd)
$editRng.Text = $new

# Shift boundaries that fall after the edit point by the length delta.
$shifted = New-Object System.Collections.Generic.List[int]
foreach ($b in $boundaries) {
    if ($b -le $runStart) {
        [void]$shifted.Add($b)
    } else {
        [void]$shifted.Add($b + $delta)
    }
}

# Add the two internal split points needed for "{{" / "REPORTING_GUARDIAN" / "}}".
[void]$shifted.Add($runStart + 2)
[void]$shifted.Add($runStart + 2 + $innerNew.Length)

$final = @($shifted | Sort-Object -Unique)

# Re-assert each original (or newly required) run boundary by toggling
# Bold off-and-back-on across the whole sub-range between consecutive
# boundaries. This forces the engine to re-split the merged run without
# changing any actual formatting.
for ($i = 0; $i -lt $final.Count - 1; $i++) {
    $a = $final[$i]
    $b = $final[$i + 1]
    if ($b -le $a) { continue }
    $sub = $d.Range($a, $b)
    $sub.Font.Bold = 1
    $sub.Font.Bold = 0
}

Write-Output "ok"
